$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 418.90475
$ws.Range("I2").Value = 114.411766
$ws.Range("K2").Value = 114.411766
$ws.Range("M2").Value = -1.411766
# row 5
$ws.Range("H5").Value = 236.26923
$ws.Range("I5").Value = 92
$ws.Range("K5").Value = 92
$ws.Range("M5").Value = 23
# row 33
$ws.Range("H33").Value = 2039.6
$ws.Range("I33").Value = 2039.6
$ws.Range("K33").Value = 2039.6
$ws.Range("M33").Value = -1810.6
# row 40
$ws.Range("H40").Value = 6813.552
$ws.Range("I40").Value = 6216.8335
$ws.Range("J40").Value = 7234.7646
$ws.Range("K40").Value = 6216.8335
$ws.Range("L40").Value = 7234.7646
$ws.Range("M40").Value = -6041.8335
$ws.Range("N40").Value = -7584.7646
# row 58
$ws.Range("H58").Value = 4602
$ws.Range("I58").Value = 2332.5833
$ws.Range("K58").Value = 6997.749899999999
$ws.Range("M58").Value = -6847.749899999999
# row 62
$ws.Range("H62").Value = 8933349
$ws.Range("I62").Value = 15627376
$ws.Range("K62").Value = 15627376
$ws.Range("M62").Value = -15626752
# row 64
$ws.Range("H64").Value = 10490.6
$ws.Range("J64").Value = 11863.25
$ws.Range("L64").Value = 11863.25
$ws.Range("N64").Value = -12359.25
# row 65
$ws.Range("H65").Value = 8933349
$ws.Range("I65").Value = 15627376
$ws.Range("K65").Value = 78136880
$ws.Range("M65").Value = -78133760
# row 67
$ws.Range("H67").Value = 10490.6
$ws.Range("J67").Value = 11863.25
$ws.Range("L67").Value = 11863.25
$ws.Range("N67").Value = -13579.25
# row 70
$ws.Range("H70").Value = 169833.33
$ws.Range("J70").Value = 203500
$ws.Range("L70").Value = 610500
$ws.Range("N70").Value = -611040
# row 73
$ws.Range("H73").Value = 169833.33
$ws.Range("J73").Value = 203500
$ws.Range("L73").Value = 610500
$ws.Range("N73").Value = -612372
# row 116
$ws.Range("H116").Value = 20225
$ws.Range("I116").Value = 6984.5
$ws.Range("J116").Value = 25521.2
$ws.Range("K116").Value = 6984.5
$ws.Range("L116").Value = 25521.2
$ws.Range("M116").Value = -3542.5
$ws.Range("N116").Value = -32405.2
# row 132
$ws.Range("H132").Value = 6609
$ws.Range("I132").Value = 7913.75
$ws.Range("K132").Value = 23741.25
$ws.Range("M132").Value = -21211.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 2869.7856
$ws.Range("I45").Value = 1840.8572
$ws.Range("J45").Value = 3898.7144
$ws.Range("K45").Value = 1840.8572
$ws.Range("L45").Value = 3898.7144
$ws.Range("M45").Value = -1463.8572
$ws.Range("N45").Value = -4652.7144

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 82
$ws.Range("H82").Value = 96666.336
# row 85
$ws.Range("H85").Value = 96666.336
# row 134
$ws.Range("H134").Value = 35804.195
$ws.Range("J134").Value = 96440.91
$ws.Range("L134").Value = 289322.73
$ws.Range("N134").Value = -294392.73

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 353.33334
$ws.Range("I7").Value = 341.2
$ws.Range("K7").Value = 341.2
$ws.Range("M7").Value = -228.2
# row 22
$ws.Range("H22").Value = 274.1
$ws.Range("I22").Value = 100.25
$ws.Range("J22").Value = 390
$ws.Range("K22").Value = 100.25
$ws.Range("L22").Value = 390
$ws.Range("M22").Value = 249.75
$ws.Range("N22").Value = -1090
# row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# row 132
$ws.Range("H132").Value = 4733.931
$ws.Range("I132").Value = 4561.9473
$ws.Range("J132").Value = 5060.7
$ws.Range("K132").Value = 13685.8419
$ws.Range("L132").Value = 15182.1
$ws.Range("M132").Value = -11155.8419
$ws.Range("N132").Value = -20242.1
# row 134
$ws.Range("H134").Value = 1056492.6
$ws.Range("I134").Value = 670423.1
$ws.Range("K134").Value = 2011269.3
$ws.Range("M134").Value = -2008734.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 34
$ws.Range("H34").Value = 102808.73
$ws.Range("J34").Value = 102808.73
$ws.Range("L34").Value = 308426.19
$ws.Range("N34").Value = -308594.19
# row 68
$ws.Range("H68").Value = 1334449.8
$ws.Range("I68").Value = 1250480.1
$ws.Range("J68").Value = 1430415.1
$ws.Range("K68").Value = 3751440.3
$ws.Range("L68").Value = 4291245.300000001
$ws.Range("M68").Value = -3750629.3
$ws.Range("N68").Value = -4292867.300000001
# row 71
$ws.Range("H71").Value = 1334449.8
$ws.Range("I71").Value = 1250480.1
$ws.Range("J71").Value = 1430415.1
$ws.Range("K71").Value = 11254320.9
$ws.Range("L71").Value = 12873735.9
$ws.Range("M71").Value = -11250264.9
$ws.Range("N71").Value = -12881847.9
# row 92
$ws.Range("H92").Value = 500758.6
$ws.Range("I92").Value = 667178.1
$ws.Range("K92").Value = 2001534.3
$ws.Range("M92").Value = -2000286.3
# row 122
$ws.Range("H122").Value = 44246.086
$ws.Range("I122").Value = 563.1667
$ws.Range("K122").Value = 5068.5003
$ws.Range("M122").Value = -2618.5003
# row 131
$ws.Range("H131").Value = 15246389
$ws.Range("I131").Value = 41792390
$ws.Range("J131").Value = 77245.64
$ws.Range("K131").Value = 125377170
$ws.Range("L131").Value = 231736.92
$ws.Range("M131").Value = -125372130
$ws.Range("N131").Value = -241816.92
# row 139
$ws.Range("H139").Value = 5329.32
$ws.Range("I139").Value = 4326
$ws.Range("J139").Value = 7909.2856
$ws.Range("K139").Value = 12978
$ws.Range("L139").Value = 23727.8568
$ws.Range("M139").Value = -7838
$ws.Range("N139").Value = -34007.8568

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 2689.5454
$ws.Range("I132").Value = 2062.6785
$ws.Range("K132").Value = 6188.0355
$ws.Range("M132").Value = -3658.0355
# row 135
$ws.Range("H135").Value = 63944.5
$ws.Range("J135").Value = 63944.5
$ws.Range("L135").Value = 63944.5
$ws.Range("N135").Value = -74084.5
# row 138
$ws.Range("H138").Value = 79950
$ws.Range("J138").Value = 79950
$ws.Range("L138").Value = 79950
$ws.Range("N138").Value = -90230

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 7999.8
$ws.Range("J62").Value = 7142.7144
$ws.Range("L62").Value = 7142.7144
$ws.Range("N62").Value = -8390.714400000001
# row 65
$ws.Range("H65").Value = 7999.8
$ws.Range("J65").Value = 7142.7144
$ws.Range("L65").Value = 35713.572
$ws.Range("N65").Value = -41953.572
# row 74
$ws.Range("H74").Value = 6822.727
$ws.Range("J74").Value = 7284.778
$ws.Range("L74").Value = 7284.778
$ws.Range("N74").Value = -9156.778
# row 77
$ws.Range("H77").Value = 6822.727
$ws.Range("J77").Value = 7284.778
$ws.Range("L77").Value = 21854.334
$ws.Range("N77").Value = -31214.334
# row 100
$ws.Range("H100").Value = 659.7083
$ws.Range("I100").Value = 674.65
$ws.Range("J100").Value = 585
$ws.Range("K100").Value = 1349.3
$ws.Range("L100").Value = 1170
$ws.Range("M100").Value = -808.3
$ws.Range("N100").Value = -2252
# row 126
$ws.Range("H126").Value = 4312.5
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
# row 132
$ws.Range("H132").Value = 39453.863
$ws.Range("I132").Value = 4383.5
$ws.Range("J132").Value = 117388
$ws.Range("K132").Value = 13150.5
$ws.Range("L132").Value = 352164
$ws.Range("M132").Value = -10620.5
$ws.Range("N132").Value = -357224
# row 136
$ws.Range("H136").Value = 8826242
$ws.Range("I136").Value = 10753876
$ws.Range("J136").Value = 289576.84
$ws.Range("K136").Value = 32261628
$ws.Range("L136").Value = 868730.52
$ws.Range("M136").Value = -32259078
# row 138
$ws.Range("H138").Value = 80346.664
$ws.Range("J138").Value = 80346.664
$ws.Range("L138").Value = 80346.664
$ws.Range("N138").Value = -90626.664
